# Add a new "Pixels per cm" setting row to the flytrack_settings sheet,
# to support heat and light rigs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: label in column A, value in column B.
$ws.Range("A19").Value = "Pixels per cm"
$ws.Range("B19").Value = 108

# Move/update the active selection to the newly added label cell.
$ws.Range("A19").Select()
